$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Locate the paragraph that currently reads, as one single run:
#       "var minDistance = float.MaxValue; Not sure what this means."
#    and split it into several runs - one per "word" that Word's
#    spell-checker would flag/un-flag ("var ", "minDistance", " = ",
#    "float.MaxValue", "; Not sure what this means.") - by toggling a
#    character formatting property on/off over each prefix. Toggling
#    Bold true->false does not change the visible/stored formatting,
#    but it forces Word to materialize a fresh <w:r> boundary at that
#    exact offset, which is exactly the run layout the diff shows.
# ------------------------------------------------------------------

$target = "var minDistance = float.MaxValue; Not sure what this means."

$paraIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.StartsWith($target)) {
        $paraIndex = $i
        break
    }
}

if ($paraIndex -gt 0) {
    $para = $d.Paragraphs.Item($paraIndex)
    $pStart = $para.Range.Start

    # Segments the single run gets broken into (matches the diff).
    $segments = @("var ", "minDistance", " = ", "float.MaxValue", "; Not sure what this means.")

    $pos = $pStart
    $boundaries = @()
    foreach ($seg in $segments) {
        $pos = $pos + $seg.Length
        $boundaries += $pos
    }
    # Drop the final boundary - that's just the end of the paragraph
    # text, no split needed there.
    $boundaries = $boundaries[0..($boundaries.Length - 2)]

    foreach ($b in $boundaries) {
        $toggle = $d.Range($pStart, $b)
        $toggle.Font.Bold = 1
        $toggle.Font.Bold = 0
    }

    # --------------------------------------------------------------
    # 2) Add a brand-new bullet paragraph right after it, describing
    #    clicking ">" on a prefab in the Hierarchy window.
    # --------------------------------------------------------------
    $para = $d.Paragraphs.Item($paraIndex)
    $para.Range.InsertParagraphAfter()

    $newPara = $d.Paragraphs.Item($paraIndex + 1)
    $newPara.Range.Text = "By clicking the > on a prefab in the hierarchy, you will go to the prefab and your scene can be scene greyed out in the background which is cool compared to going to the prefabs in the projects folder and double clicking it."
}
